$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New poll rows (ifop, included) -> rows 61-63, and (odoxa, partially) -> rows 64-66
# id, year, week, month, day, firm, collectmode, unsure, n, c_poutou, c_arthaud,
# c_melenchon, c_roussel, c_montebourg, c_jadot, c_hidalgo, c_macron, c_pecresse,
# c_barnier, c_bertrand, c_lassalle, c_daignant, c_lepen, c_zemmour

$ws.Range("A61").Value = 18
$ws.Range("B61").Value = 2021
$ws.Range("C61").Value = 7
$ws.Range("D61").Value = 10
$ws.Range("E61").Value = 11
$ws.Range("F61").Value = "ifop"
$ws.Range("G61").Value = "online"
$ws.Range("H61").Value = "included"
$ws.Range("I61").Value = 3432
$ws.Range("J61").Value = 0.5
$ws.Range("K61").Value = 0.5
$ws.Range("L61").Value = 8
$ws.Range("M61").Value = 1.5
$ws.Range("N61").Value = 1.5
$ws.Range("O61").Value = 7
$ws.Range("P61").Value = 5
$ws.Range("Q61").Value = 25
$ws.Range("T61").Value = 15
$ws.Range("U61").Value = 0.5
$ws.Range("V61").Value = 2.5
$ws.Range("W61").Value = 17
$ws.Range("X61").Value = 16

$ws.Range("A62").Value = 18
$ws.Range("B62").Value = 2021
$ws.Range("C62").Value = 7
$ws.Range("D62").Value = 10
$ws.Range("E62").Value = 11
$ws.Range("F62").Value = "ifop"
$ws.Range("G62").Value = "online"
$ws.Range("H62").Value = "included"
$ws.Range("I62").Value = 3430
$ws.Range("J62").Value = 0.5
$ws.Range("K62").Value = 0.5
$ws.Range("L62").Value = 8
$ws.Range("M62").Value = 2
$ws.Range("N62").Value = 2
$ws.Range("O62").Value = 7
$ws.Range("P62").Value = 5.5
$ws.Range("Q62").Value = 26
$ws.Range("R62").Value = 10
$ws.Range("U62").Value = 1
$ws.Range("V62").Value = 2.5
$ws.Range("W62").Value = 18
$ws.Range("X62").Value = 17

$ws.Range("A63").Value = 18
$ws.Range("B63").Value = 2021
$ws.Range("C63").Value = 7
$ws.Range("D63").Value = 10
$ws.Range("E63").Value = 11
$ws.Range("F63").Value = "ifop"
$ws.Range("G63").Value = "online"
$ws.Range("H63").Value = "included"
$ws.Range("I63").Value = 3429
$ws.Range("J63").Value = 0.5
$ws.Range("K63").Value = 0.5
$ws.Range("L63").Value = 8
$ws.Range("M63").Value = 2
$ws.Range("N63").Value = 2
$ws.Range("O63").Value = 7
$ws.Range("P63").Value = 5.5
$ws.Range("Q63").Value = 27
$ws.Range("S63").Value = 8
$ws.Range("U63").Value = 1
$ws.Range("V63").Value = 3
$ws.Range("W63").Value = 18.5
$ws.Range("X63").Value = 17

$ws.Range("A64").Value = 19
$ws.Range("B64").Value = 2021
$ws.Range("C64").Value = 6
$ws.Range("D64").Value = 10
$ws.Range("E64").Value = 8
$ws.Range("F64").Value = "odoxa"
$ws.Range("G64").Value = "online"
$ws.Range("H64").Value = "partially"
$ws.Range("I64").Value = 1856
$ws.Range("J64").Value = 2
$ws.Range("K64").Value = 1
$ws.Range("L64").Value = 8
$ws.Range("M64").Value = 2
$ws.Range("N64").Value = 2
$ws.Range("O64").Value = 6.5
$ws.Range("P64").Value = 4.5
$ws.Range("Q64").Value = 25
$ws.Range("T64").Value = 13
$ws.Range("V64").Value = 2
$ws.Range("W64").Value = 18
$ws.Range("X64").Value = 16

$ws.Range("A65").Value = 19
$ws.Range("B65").Value = 2021
$ws.Range("C65").Value = 6
$ws.Range("D65").Value = 10
$ws.Range("E65").Value = 8
$ws.Range("F65").Value = "odoxa"
$ws.Range("G65").Value = "online"
$ws.Range("H65").Value = "partially"
$ws.Range("I65").Value = 1942
$ws.Range("J65").Value = 2
$ws.Range("K65").Value = 1
$ws.Range("L65").Value = 9
$ws.Range("M65").Value = 2
$ws.Range("N65").Value = 2.5
$ws.Range("O65").Value = 7
$ws.Range("P65").Value = 4
$ws.Range("Q65").Value = 26.5
$ws.Range("R65").Value = 6.5
$ws.Range("V65").Value = 2.5
$ws.Range("W65").Value = 18.5
$ws.Range("X65").Value = 16.5

$ws.Range("A66").Value = 19
$ws.Range("B66").Value = 2021
$ws.Range("C66").Value = 6
$ws.Range("D66").Value = 10
$ws.Range("E66").Value = 8
$ws.Range("F66").Value = "odoxa"
$ws.Range("G66").Value = "online"
$ws.Range("H66").Value = "partially"
$ws.Range("I66").Value = 1917
$ws.Range("J66").Value = 2
$ws.Range("K66").Value = 1
$ws.Range("L66").Value = 8.5
$ws.Range("M66").Value = 2
$ws.Range("O66").Value = 11
$ws.Range("Q66").Value = 26
$ws.Range("T66").Value = 13
$ws.Range("V66").Value = 2.5
$ws.Range("W66").Value = 18
$ws.Range("X66").Value = 16

# Update view state to match: scrolled so column F is leftmost in the top pane,
# frozen pane scrolled to row 58, and selection on the bottom-right-most new cell.
$win = $excel.ActiveWindow
$win.ScrollColumn = 6
$win.ScrollRow = 58
$ws.Range("Z66").Select() | Out-Null
